$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data edit: altitude reading for KAL85 @ DCA (F7) revised from 28000 to 27500
$ws.Range("F7").Value = 27500

# Update the active cell/selection to reflect where the editor left off (F8)
$ws.Range("F8").Select()
